# Add a new bug row (row 4) to the Bug Report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Bug_Client_001"
$ws.Range("B4").Value = "TC_Client_017"
$ws.Range("C4").Value = "Empty cart message not appear"
$ws.Range("D4").Value = "no message appear"
$ws.Range("E4").Value = " message: ""Your cart is empty"""
$ws.Range("F4").Value = "1. Login with Client email`n2.  Navigate to Home page `n3. Click on Cart Button`n4. Message will appear"
$ws.Range("F4").WrapText = $true
$ws.Range("G4").Value = "Medium"
$ws.Range("H4").Value = "Medium"
$ws.Range("I4").Value = "Development"
$ws.Range("J4").Value = "Abdallah"
$ws.Range("K4").Value = "Closed"

# Move the active selection to the newly added row, matching the saved
# workbook's cursor position.
$ws.Range("K4").Select()
